$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New "Notes" column header (D1) ---
$ws.Range("D1").Value = "Notes"

# --- Edit existing rows 2-4 (source/amount/date revised; row4 gets a new amount) ---
# C2/C3/C4 already carry the workbook's original date style (s="1"),
# so only their values need to change - no NumberFormat touch needed there.
$ws.Range("A2").Value = "pppp"
$ws.Range("B2").Value = 5000
$ws.Range("C2").Value = 45982.23976851852
$ws.Range("D2").Value = "helooiii"

$ws.Range("A3").Value = "Salary"
$ws.Range("B3").Value = 25000
$ws.Range("C3").Value = 45966.23976851852
$ws.Range("D3").Value = ""

$ws.Range("A4").Value = "Salary"
$ws.Range("B4").Value = 5000
$ws.Range("C4").Value = 45966.23976851852
$ws.Range("D4").Value = ""

# --- Append new rows 5-13 (new income/expense entries with the Notes column) ---
# Each row: Source, Amount, Date (serial), Notes
$newRows = @(
    @{ Row = 5;  Source = "Salary2";  Amount = 50000;  Date = 45966.23976851852; Notes = "" },
    @{ Row = 6;  Source = "random";   Amount = 234345; Date = 45966.23976851852; Notes = "yoyoyo" },
    @{ Row = 7;  Source = "fsdfvds";  Amount = 12345;  Date = 45926.23976851852; Notes = "" },
    @{ Row = 8;  Source = "salary2";  Amount = 25367;  Date = 45925.23976851852; Notes = "" },
    @{ Row = 9;  Source = "Interest"; Amount = 50000;  Date = 45922.23976851852; Notes = "" },
    @{ Row = 10; Source = "utilities";Amount = 2386;   Date = 45905.23976851852; Notes = "" },
    @{ Row = 11; Source = "salary";   Amount = 50000;  Date = 45904.23976851852; Notes = "" },
    @{ Row = 12; Source = "salary1";  Amount = 50000;  Date = 45904.23976851852; Notes = "" },
    @{ Row = 13; Source = "salary";   Amount = 50000;  Date = 45870.23976851852; Notes = "" }
)

foreach ($r in $newRows) {
    $ws.Range("A$($r.Row)").Value = $r.Source
    $ws.Range("B$($r.Row)").Value = $r.Amount
    $ws.Range("C$($r.Row)").Value = $r.Date
    # Brand-new cells have no inherited style, so stamp them with the same
    # short-date format ("m/d/yy" == built-in numFmtId 14) used by C2:C4.
    $ws.Range("C$($r.Row)").NumberFormat = "m/d/yy"
    $ws.Range("D$($r.Row)").Value = $r.Notes
}
